$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Variables")

# child_id's valueType is converted back to "text".
$ws.Cells.Item(2, 2).Value = "text"

# Insert a new row 2 for the "row_id" variable (shifts existing rows down by one,
# so child_id - just edited above - becomes row 3).
$ws.Rows.Item(2).Insert()

# Populate the new row_id row.
$ws.Cells.Item(2, 1).Value = "row_id"
$ws.Cells.Item(2, 2).Value = "integer"
$ws.Cells.Item(2, 3).Value = "numeric"
$ws.Cells.Item(2, 4).Value = "Unique identifier for the row in Opal"

# Format the new row: text number format plus a thin white border, to mark it
# as a distinct/new entry, matching the styling used elsewhere in the sheet.
$newRowRange = $ws.Range("A2:D2")
$newRowRange.NumberFormat = "@"
$newRowRange.Borders.ColorIndex = 2
$newRowRange.Borders.Weight = 2
$newRowRange.Borders.LineStyle = 1

# Highlight/select the newly inserted row, as left by the editor.
$newRowRange.Select()
